$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C31) from 2024-10-24 (45589)
# to 2024-10-25 (45590) for all data rows.
$ws.Range("C2:C31").Value = 45590
